$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("004-PI588478_15-3", 0.9857791880580289),
    @("001-PI588568_13-107", 0.8493184792593224),
    @("009-PI588470_15-17", 0.9720587635652564),
    @("012-PI588471_15-23", 0.7082437918797061),
    @("007-PI588601_15-11", 0.6993877289623024),
    @("006-PI588606_15-9", 0.9926110729154066),
    @("011-PI588463_15-21", 0.8298248193246827),
    @("010-PI588543_15-19", 0.9614144369530073),
    @("002-PI588331_13-109", 0.9649735094556326),
    @("003-PI597269_15-1", 0.8681984001287604),
    @("008-PI588487_15-15", 0.8292395127984116),
    @("005-PI588133_15-5", 0.9903248612495549)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
